# Weekend catch-up / tbd empty row fix
# Appends 4 new contract rows (106-109) to the WIP tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 106; A = "4/25/2019"; B = "SPE5EM-19-V-4186"; C = "9";  D = '$324.00 ';    E = "5330012549189"; F = "GASKET";              G = "American Metal Bearing"; H = "G001";       I = "M33"; J = "2019 OCT 02" },
    @{ Row = 107; A = "4/26/2019"; B = "SPE4A6-19-P-E296"; C = "6";  D = '$13,695.42';   E = "5895016527348"; F = "TRANSMITTER GROUP";   G = "GEMS";                    H = "822610";     I = "M41"; J = "2019 OCT 07" },
    @{ Row = 108; A = "4/26/2019"; B = "SPE7M8-19-P-2323"; C = "7";  D = '$16,436.00';   E = "5930014065242"; F = "SWITCH,LIQUID LEVEL"; G = "GEMS";                    H = "LS-76725";   I = "CP";  J = "2019 OCT 11" },
    @{ Row = 109; A = "4/26/2019"; B = "SPE4A6-19-V-136J"; C = "10"; D = '$61,583.70 ';  E = "6680013650925"; F = "TRANSMITTER,LIQUID";  G = "GEMS";                    H = "42880-0105"; I = "CP";  J = "2019 OCT 08" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $rowRange = $ws.Range("A" + $rowNum + ":J" + $rowNum)

    # Force every cell in the row to be stored as text, matching the rest of
    # the sheet (all data cells are plain strings, even the numeric-looking
    # ones like Qty or NSN). Without this, values like "9" or "4/25/2019"
    # would be auto-coerced to numbers/dates by the COM value setter.
    $rowRange.NumberFormat = "@"

    $ws.Range("A" + $rowNum).Value = $r.A
    $ws.Range("B" + $rowNum).Value = $r.B
    $ws.Range("C" + $rowNum).Value = $r.C
    $ws.Range("D" + $rowNum).Value = $r.D
    $ws.Range("E" + $rowNum).Value = $r.E
    $ws.Range("F" + $rowNum).Value = $r.F
    $ws.Range("G" + $rowNum).Value = $r.G
    $ws.Range("H" + $rowNum).Value = $r.H
    $ws.Range("I" + $rowNum).Value = $r.I
    $ws.Range("J" + $rowNum).Value = $r.J

    # Drop the "@" text number-format again so the saved cells carry no
    # style override (matches the rest of the sheet, which has no `s`
    # attribute on its text cells) while keeping the text cell type.
    $rowRange.ClearFormats()
}
